$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 227 (existing rows 227-338 shift down to 230-341)
$ws.Rows("227:229").Insert()

# Values that are constant across the whole dataset (same market/product)
$constA = 3
$constB = "Femacal de La Calera"
$constC = "Coquimbo"
$constE = 5
$constF = "Fruta"
$constG = 100101
$constH = "Berries"
$constI = 100112025
$constJ = "Frutilla"
$constK = "Sin especificar"
$constQ = '$/bandeja 7 kilos'
$constR = "Provincia de Melipilla"
$constT = 7

$newRows = @(
    @{ Row = 227; D = 44839; L = "Especial"; M = 45; N = 13000; O = 13000; P = 13000; S = 1857 },
    @{ Row = 228; D = 44839; L = "Primera";  M = 50; N = 10000; O = 10000; P = 10000; S = 1429 },
    @{ Row = 229; D = 44839; L = "Segunda";  M = 36; N = 8000;  O = 8000;  P = 8000;  S = 1143 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $constA
    $ws.Cells.Item($row, 2).Value = $constB
    $ws.Cells.Item($row, 3).Value = $constC
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $constE
    $ws.Cells.Item($row, 6).Value = $constF
    $ws.Cells.Item($row, 7).Value = $constG
    $ws.Cells.Item($row, 8).Value = $constH
    $ws.Cells.Item($row, 9).Value = $constI
    $ws.Cells.Item($row, 10).Value = $constJ
    $ws.Cells.Item($row, 11).Value = $constK
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $constQ
    $ws.Cells.Item($row, 18).Value = $constR
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = $constT
}
